# Regenerate orders with updated distance/size codes.
# Distances: D64 -> D69, D51 -> D55, D80 -> D86
# Sizes:     S30 -> S31
# These substrings are embedded inside many composite condition / filename
# strings (e.g. "Face02_D64_S25", "Face18_D64_S30_l.png") as well as
# appearing standalone in the Distance/Size lookup columns, so a workbook
# wide substring Find&Replace reproduces every edit in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

$cells.Replace("D64", "D69")
$cells.Replace("D51", "D55")
$cells.Replace("D80", "D86")
$cells.Replace("S30", "S31")
